$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "New Retailers Format" -- update/add retailer rows 2-4
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 2 - Al Faruk Electronics (replaces the previous Utsob Telecom entry)
$ws1.Range("A2").Value = "DEL-0179"
$ws1.Range("B2").Value = "DSR-0247"
$ws1.Range("C2").Value = "Al Faruk Electronics"
$ws1.Range("D2").Value = "Naldanga"
$ws1.Range("E2").Value = "Md Faruk Hosen Mukta"
$ws1.Range("G2").Value = "GO"
$ws1.Range("I2").Value = "Md Faruk Hosen Mukta"
$ws1.Range("J2").Value = 1733193387
$ws1.Range("K2").Value = "Natore"
$ws1.Range("L2").Value = "Natore Sodor"
$ws1.Range("M2").Value = "ZSO-0022"
$ws1.Range("N2").Value = "3/1 RB Super Market, Station Bazar, Natore."
$ws1.Range("P2").Value = 1733193387
$ws1.Range("Q2").Value = "C"
$ws1.Range("R2").Value = "Rural"
$ws1.Range("S2").Value = "bKash"
$ws1.Range("T2").Value = 1733193387

# Row 3 - FA Phone
$ws1.Range("A3").Value = "DEL-0179"
$ws1.Range("B3").Value = "DSR-0248"
$ws1.Range("C3").Value = "FA Phone"
$ws1.Range("D3").Value = "Baraigram"
$ws1.Range("E3").Value = "FA Shohel"
$ws1.Range("G3").Value = "GO"
$ws1.Range("I3").Value = "FA Shohel"
$ws1.Range("J3").Value = 1717299513
$ws1.Range("K3").Value = "Natore"
$ws1.Range("L3").Value = "Baraigram"
$ws1.Range("M3").Value = "ZSO-0022"
$ws1.Range("N3").Value = "Jonail Bazar, Baraigram, Natore."
$ws1.Range("P3").Value = 1717299513
$ws1.Range("Q3").Value = "C"
$ws1.Range("R3").Value = "Rural"
$ws1.Range("S3").Value = "bKash"
$ws1.Range("T3").Value = 1717299513

# Row 4 - Dolon Mobile
$ws1.Range("A4").Value = "DEL-0179"
$ws1.Range("B4").Value = "DSR-0248"
$ws1.Range("C4").Value = "Dolon Mobile"
$ws1.Range("D4").Value = "Baraigram"
$ws1.Range("E4").Value = "Dolon"
$ws1.Range("G4").Value = "GO"
$ws1.Range("I4").Value = "Dolon"
$ws1.Range("K4").Value = "Natore"
$ws1.Range("L4").Value = "Baraigram"
$ws1.Range("M4").Value = "ZSO-0022"
$ws1.Range("N4").Value = "Laxmicole, Baraigram, Natore."
$ws1.Range("Q4").Value = "C"
$ws1.Range("R4").Value = "Rural"
$ws1.Range("S4").Value = "bKash"

# The new, longer "Thana" (col L) and "Address" (col N) values no longer fit
# the old best-fit column widths, so widen those two columns to match.
$ws1.Columns.Item(12).ColumnWidth = 11.877604166666666
$ws1.Columns.Item(14).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# Update the "current selection" (cursor position) on each sheet so that the
# saved workbook reopens with the same cell focus as in the authored edit.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Select()
$ws2.Range("G17").Select()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Select()
$ws3.Range("O12").Select()

# Re-activate sheet 1 last so it remains the tab shown when the file opens,
# and leave the selection on the E16:E18 range (matching the authored edit).
$ws1.Select()
$ws1.Range("E16:E18").Select()
